# Update the EPEX spot prices workbook with the latest day of data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column (CU) with header "20-sep"
# and the hourly prices for that day.
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (bold / border / centered)
# onto the new header cell so it reuses the existing header style instead of
# creating a brand-new one.
$wsSpot.Range("CT1").Copy()
$wsSpot.Range("CU1").PasteSpecial(-4122)   # xlPasteFormats
$wsSpot.Range("CU1").Value = "20-sep"

$spotValues = @{
  2  = 0
  3  = 14.71
  4  = 20.02
  5  = 19.31
  6  = 20.41
  7  = 19.31
  8  = 25.64
  9  = 20.8
  10 = 18.47
  11 = 24.87
  12 = 3.39
  13 = 3.94
  14 = 1.72
  15 = 0
  16 = 0
  17 = 0.65
  18 = 4.31
  19 = 15.43
  20 = 32.52
  21 = 44.62
  22 = 68.65000000000001
  23 = 58.21
  24 = 45.32
  25 = 17.62
}

foreach ($row in $spotValues.Keys) {
  $wsSpot.Range("CU$row").Value = $spotValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append a new row with the latest date / price.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$gazDate = $wsGaz.Range("A96")
$gazDate.NumberFormat = "@"        # force text so "2025-09-18" is not
$gazDate.Value = "2025-09-18"      # auto-converted into a date serial
$gazDate.Style = "Normal"          # drop the temporary text format again

$wsGaz.Range("B96").Value = 31.8

# ---------------------------------------------------------------------
# Sheet "CO2": append a new row with the latest date / price.
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$co2Date = $wsCO2.Range("A96")
$co2Date.NumberFormat = "@"
$co2Date.Value = "2025-09-18"
$co2Date.Style = "Normal"

$wsCO2.Range("B96").Value = 77.2
